# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/E are always non-numeric-looking text (names, URLs, "  +x.xx%  ").
# Column D sometimes looks like a plain number (e.g. "245.15") - without forcing
# a text format first, Excel would auto-convert it to a floating point number
# (losing the exact decimal text, e.g. "0.661" -> 0.66100000000000003). The
# source data keeps these as literal text, so mark such cells as Text before 
# writing the value whenever the new string parses as a number.

# Row 2
$ws.Range("D2").Value = '36.447.37'
$ws.Range("E2").Value = '  +2.67%  '

# Row 3
$ws.Range("D3").Value = '2.008.37'
$ws.Range("E3").Value = '  +5.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.15'
$ws.Range("E5").Value = '  -0.53%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.661'
$ws.Range("E6").Value = '  -4.75%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.94'
$ws.Range("E8").Value = '  +4.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.38'
$ws.Range("E9").Value = '  +9.29%  '

# Row 10
$ws.Range("E10").Value = '  +3.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0714'
$ws.Range("E11").Value = '  -5.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0983'
$ws.Range("E12").Value = '  -0.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.58'
$ws.Range("E13").Value = '  +2.10%  '

# Row 14
$ws.Range("D14").Value = '2.301.96'
$ws.Range("E14").Value = '  +5.92%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.814'
$ws.Range("E15").Value = '  +2.30%  '

# Row 16
$ws.Range("D16").Value = '2.012.18'
$ws.Range("E16").Value = '  +5.85%  '

# Row 17
$ws.Range("E17").Value = '  -2.41%  '

# Row 18
$ws.Range("D18").Value = '36.325.64'
$ws.Range("E18").Value = '  +2.47%  '

# Row 19
$ws.Range("E19").Value = '  -3.26%  '

# Row 20
$ws.Range("E20").Value = '  -1.93%  '

# Row 21
$ws.Range("E21").Value = '  -1.24%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.89'
$ws.Range("E22").Value = '  -3.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.88'
$ws.Range("E23").Value = '  -6.15%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("E25").Value = '  -10.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.74'
$ws.Range("E26").Value = '  -1.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.63'
$ws.Range("E27").Value = '  -0.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.66'

# Row 29
$ws.Range("E29").Value = '  -10.84%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  -5.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.31'
$ws.Range("E31").Value = '  +62.86%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.39'
$ws.Range("E32").Value = '  +0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0585'
$ws.Range("E33").Value = '  -2.94%  '

# Row 34
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("E35").Value = '  -0.80%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.98'
$ws.Range("E36").Value = '  -6.23%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.16'
$ws.Range("E37").Value = '  +10.81%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0803'
$ws.Range("E38").Value = '  +8.84%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.855'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("E40").Value = '  -9.36%  '

# Row 41
$ws.Range("E41").Value = '  -3.85%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.16'
$ws.Range("E42").Value = '  -3.09%  '

# Row 43
$ws.Range("E43").Value = '  +2.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  +16.70%  '

# Row 45
$ws.Range("E45").Value = '  -5.82%  '

# Row 46
$ws.Range("D46").Value = '1.315.48'
$ws.Range("E46").Value = '  -0.73%  '

# Row 47
$ws.Range("E47").Value = '  +0.70%  '

# Row 48
$ws.Range("E48").Value = '  +0.79%  '

# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.201.70'
$ws.Range("E49").Value = '  +6.18%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.20'
$ws.Range("E50").Value = '  -7.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.85'
$ws.Range("E51").Value = '  +15.12%  '
